# employeedata.xlsx edit:
#  - change selection/active view on "negativelogindata" (sheet 1)
#  - append a new "personnalDetails" worksheet with a small personnel table
#    and make it the active/selected sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. negativelogindata: it used to be tabSelected with topLeftCell=A4 and a
#    B17:B18 selection; the new view just leaves B14 selected (no special
#    scroll position / tab-selected flag, because the new sheet takes that
#    spot).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("B14").Select()

# ---------------------------------------------------------------------------
# 2. Add the new "personnalDetails" sheet after the last existing sheet.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "personnalDetails"

# Header row + 4 rows of sample personnel data (12 columns: A..L)
$tableData = @(
  @("License Expiry Date", "driver licence number", "Ssn number", "Other id", "Nick name", "Military service", "Sin number", "Marital Status", "Nationality", "Smoker", "dat of birth", "gender"),
  @("1994-Mar-15", "A66666325", 788, 11111, "L", "Yes", 77, "Single", "Algerian", "yes", "1994-Mar-15", "male"),
  @("1994-Mar-15", "B44445558", 7525, 22225, "A", "No", 12, "Married ", "American", "no", "1994-Mar-15", "female"),
  @("1994-Mar-15", "J448885285", 25488, 99995, "F", "Yes", 0, "Other", "Armenian", "yes", "1994-Mar-15", "male"),
  @("1994-Mar-15", "Y4885626555", 44525, 55555, "B", "No", 56, "Single", "Bahamian", "no", "1994-Mar-15", "male")
)

for ($r = 0; $r -lt $tableData.Count; $r++) {
  $rowValues = $tableData[$r]
  for ($c = 0; $c -lt $rowValues.Count; $c++) {
    $ws4.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
  }
}

# Column widths (character units, closest values this host can reproduce)
$ws4.Columns.Item(1).ColumnWidth = 22.166666666666668
$ws4.Columns.Item(2).ColumnWidth = 21.5
$ws4.Columns.Item(3).ColumnWidth = 25.5
$ws4.Columns.Item(4).ColumnWidth = 19.0
$ws4.Columns.Item(5).ColumnWidth = 16.166666666666668
$ws4.Columns.Item(6).ColumnWidth = 15.666666666666666
$ws4.Columns.Item(7).ColumnWidth = 18.0
$ws4.Columns.Item(8).ColumnWidth = 15.333333333333334
$ws4.Columns.Item(9).ColumnWidth = 14.166666666666666
$ws4.Columns.Item(11).ColumnWidth = 14.833333333333334

# Make the new sheet the active / selected tab, with J1 selected.
[void]$ws4.Activate()
[void]$ws4.Range("J1").Select()
